# uurschema_voorbeeld.xlsx - "Add files via upload" edit
#
# The original sheet has 4 columns (dag, ingetikt, uitgetikt, wacht) and
# 6 data rows (2021-10-01 .. 2021-10-06, then 2021-11-11). The "wacht"
# column only had meaningful (TRUE) values for 2021-10-03 and 2021-11-11;
# the FALSE placeholders are removed entirely.
#
# Two new rows are inserted for 2021-10-07 and 2021-10-08 (before the
# 2021-11-11 row), and two new boolean columns are added: "verlof" (E)
# and "recup" (F). 2021-10-07 is marked recup=TRUE, 2021-10-08 is marked
# verlof=TRUE. The newly inserted rows only have a day value (A) - the
# ingetikt/uitgetikt (B/C) cells stay empty but keep the date/time style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 7 (pushes the 2021-11-11 row from row 7 to
# row 9); the inserted rows inherit the formatting of the row they push
# down, giving B7/C7/B8/C8 the same date-time style as the rest of column
# B/C.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# Drop the FALSE "wacht" placeholders on the days that are not on-call;
# ClearContents() removes the cell entirely (it has no explicit style),
# matching the diff which deletes the <c r="D#" t="b"> elements outright.
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
# (D4, the 2021-10-03 wacht=TRUE value, and the old D7/new D9 2021-11-11
# wacht=TRUE value are both left untouched.)

# New header columns.
$ws.Cells.Item(1, 5).Value = "verlof"
$ws.Cells.Item(1, 6).Value = "recup"

# New row for 2021-10-07 -> recup = TRUE.
$ws.Cells.Item(7, 1).Value = 44476
$ws.Cells.Item(7, 6).Value = $true

# New row for 2021-10-08 -> verlof = TRUE.
$ws.Cells.Item(8, 1).Value = 44477
$ws.Cells.Item(8, 5).Value = $true

# Match the saved selection from the diff (D5:D7, active cell D5).
$ws.Range("D5:D7").Select() | Out-Null
